$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2022-02-28"
$ws.Range("I2").Value = "1a (cosecha)"
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 480
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = 490
$ws.Range("O2").Value = "Región de O'Higgins"
$ws.Range("P2").Value = 490

# Row 3
$ws.Range("D3").Value = "2021-06-02"
$ws.Range("H3").Value = "Pachia"
$ws.Range("I3").Value = "1a nueva(o)"
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 730
$ws.Range("L3").Value = 750
$ws.Range("M3").Value = 740
$ws.Range("P3").Value = 740

# Row 4
$ws.Range("D4").Value = "2022-04-14"
$ws.Range("I4").Value = "1a (cosecha)"
$ws.Range("J4").Value = 1200
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 420
$ws.Range("M4").Value = 410
$ws.Range("O4").Value = "Región de O'Higgins"
$ws.Range("P4").Value = 410

# Row 5
$ws.Range("D5").Value = "2022-01-14"
$ws.Range("I5").Value = "1a nueva(o)"
$ws.Range("J5").Value = 1300
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 550
$ws.Range("M5").Value = 525
$ws.Range("O5").Value = "Región de O'Higgins"
$ws.Range("P5").Value = 525

# Row 6
$ws.Range("D6").Value = "2021-12-17"
$ws.Range("I6").Value = "1a (cosecha)"
$ws.Range("J6").Value = 800

# Row 7
$ws.Range("D7").Value = "2021-12-17"
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 550
$ws.Range("M7").Value = 525
$ws.Range("O7").Value = "Perú"
$ws.Range("P7").Value = 525

# Row 8
$ws.Range("D8").Value = "2021-05-06"
$ws.Range("I8").Value = "1a (cosecha)"
$ws.Range("J8").Value = 1200
$ws.Range("K8").Value = 350
$ws.Range("L8").Value = 400
$ws.Range("M8").Value = 375
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 375

# Row 9
$ws.Range("D9").Value = "2022-05-04"
$ws.Range("I9").Value = "1a (cosecha)"
$ws.Range("K9").Value = 680
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = 690
$ws.Range("O9").Value = "Región de O'Higgins"
$ws.Range("P9").Value = 690

# Row 10
$ws.Range("D10").Value = "2021-10-14"
$ws.Range("I10").Value = "1a nueva(o)"
$ws.Range("J10").Value = 1300
$ws.Range("K10").Value = 550
$ws.Range("L10").Value = 580
$ws.Range("M10").Value = 565
$ws.Range("O10").Value = "Perú"
$ws.Range("P10").Value = 565

# Row 11
$ws.Range("D11").Value = "2020-12-10"
$ws.Range("I11").Value = "1a nueva(o)"
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1550
$ws.Range("M11").Value = 1525
$ws.Range("O11").Value = "Perú"
$ws.Range("P11").Value = 1525

# Row 12
$ws.Range("D12").Value = "2022-03-30"
$ws.Range("I12").Value = "2a (cosecha)"
$ws.Range("J12").Value = 1300
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 430
$ws.Range("M12").Value = 415
$ws.Range("O12").Value = "Región de O'Higgins"
$ws.Range("P12").Value = 415

# Row 13
$ws.Range("D13").Value = "2021-09-30"
$ws.Range("I13").Value = "1a nueva(o)"
$ws.Range("J13").Value = 1200
$ws.Range("K13").Value = 600
$ws.Range("L13").Value = 650
$ws.Range("M13").Value = 625
$ws.Range("O13").Value = "Perú"
$ws.Range("P13").Value = 625

# Row 14
$ws.Range("D14").Value = "2021-04-15"
$ws.Range("I14").Value = "2a nueva(o)"
$ws.Range("J14").Value = 1200
$ws.Range("K14").Value = 400
$ws.Range("L14").Value = 430
$ws.Range("M14").Value = 415
$ws.Range("O14").Value = "Provincia de Melipilla"
$ws.Range("P14").Value = 415

# Row 15
$ws.Range("D15").Value = "2021-02-11"
$ws.Range("J15").Value = 1250
$ws.Range("K15").Value = 430
$ws.Range("L15").Value = 450
$ws.Range("M15").Value = 440
$ws.Range("P15").Value = 440

# Row 16
$ws.Range("D16").Value = "2021-07-22"
$ws.Range("I16").Value = "1a (guarda)"
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 450
$ws.Range("L16").Value = 480
$ws.Range("M16").Value = 465
$ws.Range("O16").Value = "Provincia de Melipilla"
$ws.Range("P16").Value = 465

# Row 17
$ws.Range("D17").Value = "2021-02-16"
$ws.Range("I17").Value = "2a nueva(o)"
$ws.Range("K17").Value = 450
$ws.Range("L17").Value = 480
$ws.Range("M17").Value = 465
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 465

# Row 18
$ws.Range("D18").Value = "2021-08-20"
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 580
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = 590
$ws.Range("P18").Value = 590

# Row 19
$ws.Range("D19").Value = "2021-11-30"
$ws.Range("I19").Value = "2a nueva(o)"
$ws.Range("J19").Value = 900
$ws.Range("K19").Value = 480
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = 490
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 490

# Row 20
$ws.Range("D20").Value = "2021-01-08"
$ws.Range("I20").Value = "2a nueva(o)"
$ws.Range("J20").Value = 1600
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 550
$ws.Range("M20").Value = 525
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 525

# Row 21
$ws.Range("D21").Value = "2022-04-19"
$ws.Range("I21").Value = "1a (cosecha)"
$ws.Range("J21").Value = 1200
$ws.Range("K21").Value = 400
$ws.Range("L21").Value = 430
$ws.Range("M21").Value = 415
$ws.Range("O21").Value = "Región de O'Higgins"
$ws.Range("P21").Value = 415

# Row 22
$ws.Range("D22").Value = "2021-09-14"
$ws.Range("I22").Value = "1a nueva(o)"
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 630
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = 640
$ws.Range("O22").Value = "Perú"
$ws.Range("P22").Value = 640

# Row 23
$ws.Range("D23").Value = "2022-02-09"
$ws.Range("H23").Value = "Camote"
$ws.Range("I23").Value = "2a (cosecha)"
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 400
$ws.Range("L23").Value = 450
$ws.Range("M23").Value = 425
$ws.Range("O23").Value = "Región de O'Higgins"
$ws.Range("P23").Value = 425

# Row 24
$ws.Range("D24").Value = "2021-02-04"
$ws.Range("K24").Value = 450
$ws.Range("L24").Value = 480
$ws.Range("M24").Value = 465
$ws.Range("O24").Value = "Perú"
$ws.Range("P24").Value = 465

# Row 25
$ws.Range("D25").Value = "2021-01-15"
$ws.Range("J25").Value = 1600
$ws.Range("K25").Value = 500
$ws.Range("L25").Value = 550
$ws.Range("M25").Value = 525
$ws.Range("O25").Value = "Región de O'Higgins"
$ws.Range("P25").Value = 525

# Row 26
$ws.Range("D26").Value = "2021-10-21"
$ws.Range("I26").Value = "1a nueva(o)"
$ws.Range("O26").Value = "Perú"

# Row 27
$ws.Range("D27").Value = "2022-03-10"
$ws.Range("I27").Value = "1a (cosecha)"
$ws.Range("O27").Value = "Región Metropolitana"

# Row 28
$ws.Range("D28").Value = "2021-10-07"
$ws.Range("I28").Value = "1a nueva(o)"
$ws.Range("J28").Value = 1200
$ws.Range("K28").Value = 480
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 490
$ws.Range("O28").Value = "Perú"
$ws.Range("P28").Value = 490

# Row 29
$ws.Range("D29").Value = "2022-05-10"
$ws.Range("I29").Value = "1a (cosecha)"
$ws.Range("J29").Value = 700
$ws.Range("K29").Value = 580
$ws.Range("L29").Value = 600
$ws.Range("M29").Value = 590
$ws.Range("O29").Value = "Región de O'Higgins"
$ws.Range("P29").Value = 590

# Row 30
$ws.Range("D30").Value = "2020-12-14"
$ws.Range("I30").Value = "1a nueva(o)"
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 1350
$ws.Range("L30").Value = 1400
$ws.Range("M30").Value = 1375
$ws.Range("O30").Value = "Perú"
$ws.Range("P30").Value = 1375

# Row 31
$ws.Range("D31").Value = "2021-01-05"
$ws.Range("I31").Value = "1a nueva(o)"
$ws.Range("J31").Value = 1360
$ws.Range("K31").Value = 730
$ws.Range("L31").Value = 750
$ws.Range("M31").Value = 740
$ws.Range("O31").Value = "Perú"
$ws.Range("P31").Value = 740
